$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previously existing content (keeps row1 header cell styling),
# we will rewrite everything so the shared-string table is rebuilt in the
# exact order required by the target file.
$ws.Range("A1:T7").ClearContents()

# --- Row 1: headers (unchanged, restores shared-string indices 0..19) ---
$headers = @(
  "Sending cluster",
  "Ligand symbol",
  "Receptor symbol",
  "Target cluster",
  "Ligand-expressing cells",
  "Ligand detection rate",
  "Ligand average expression value",
  "Ligand total expression value",
  "Ligand derived specificity of average expression value",
  "Ligand derived specificity of total expression value",
  "Receptor-expressing cells",
  "Receptor detection rate",
  "Receptor average expression value",
  "Receptor total expression value",
  "Receptor derived specificity of average expression value",
  "Receptor derived specificity of total expression value",
  "Edge average expression weight",
  "Edge total expression weight",
  "Edge average expression derived specificity",
  "Edge total expression derived specificity"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
  $col = [char](65 + $i)
  $ws.Range("$col`1").Value = $headers[$i]
}

# Seed the brand-new shared strings ("sCs" then "ECs") before writing the
# rest of the grid, so they are registered in that relative order.
$ws.Range("D4").Value = "sCs"
$ws.Range("D2").Value = "ECs"

# --- Row 2 ---
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cxcl13"
$ws.Range("C2").Value = "Ackr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.009847666666666
$ws.Range("H2").Value = 12.029543
$ws.Range("I2").Value = 0.9697248931871538
$ws.Range("J2").Value = 0.9697248931871538
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1143813333333333
$ws.Range("N2").Value = 0.343144
$ws.Range("O2").Value = 0.03414257747883775
$ws.Range("P2").Value = 0.03414257747883775
$ws.Range("Q2").Value = 0.4586517225768889
$ws.Range("R2").Value = 4.127865503192
$ws.Range("S2").Value = 0.03310890729880006
$ws.Range("T2").Value = 0.03310890729880006

# --- Row 3 ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cxcl13"
$ws.Range("C3").Value = "Ackr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.009847666666666
$ws.Range("H3").Value = 12.029543
$ws.Range("I3").Value = 0.9697248931871538
$ws.Range("J3").Value = 0.9697248931871538
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.175982333333334
$ws.Range("N3").Value = 9.527947000000001
$ws.Range("O3").Value = 0.9480237703755849
$ws.Range("P3").Value = 0.9480237703755849
$ws.Range("Q3").Value = 12.73520534869122
$ws.Range("R3").Value = 114.616848138221
$ws.Range("S3").Value = 0.919322249466347
$ws.Range("T3").Value = 0.919322249466347

# --- Row 4 ---
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cxcl13"
$ws.Range("C4").Value = "Ackr4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.009847666666666
$ws.Range("H4").Value = 12.029543
$ws.Range("I4").Value = 0.9697248931871538
$ws.Range("J4").Value = 0.9697248931871538
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.05974466666666667
$ws.Range("N4").Value = 0.179234
$ws.Range("O4").Value = 0.01783365214557738
$ws.Range("P4").Value = 0.01783365214557738
$ws.Range("Q4").Value = 0.2395670122291111
$ws.Range("R4").Value = 2.156103110062
$ws.Range("S4").Value = 0.01729373642200688
$ws.Range("T4").Value = 0.01729373642200688

# --- Row 5 ---
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Cxcl13"
$ws.Range("C5").Value = "Ackr4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1251886666666667
$ws.Range("H5").Value = 0.375566
$ws.Range("I5").Value = 0.03027510681284623
$ws.Range("J5").Value = 0.03027510681284622
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1143813333333333
$ws.Range("N5").Value = 0.343144
$ws.Range("O5").Value = 0.03414257747883775
$ws.Range("P5").Value = 0.03414257747883775
$ws.Range("Q5").Value = 0.01431924661155556
$ws.Range("R5").Value = 0.128873219504
$ws.Range("S5").Value = 0.001033670180037691
$ws.Range("T5").Value = 0.001033670180037691

# --- Row 6 ---
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Cxcl13"
$ws.Range("C6").Value = "Ackr4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1251886666666667
$ws.Range("H6").Value = 0.375566
$ws.Range("I6").Value = 0.03027510681284623
$ws.Range("J6").Value = 0.03027510681284622
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.175982333333334
$ws.Range("N6").Value = 9.527947000000001
$ws.Range("O6").Value = 0.9480237703755849
$ws.Range("P6").Value = 0.9480237703755849
$ws.Range("Q6").Value = 0.3975969936668889
$ws.Range("R6").Value = 3.578372943002
$ws.Range("S6").Value = 0.02870152090923804
$ws.Range("T6").Value = 0.02870152090923804

# --- Row 7 ---
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Cxcl13"
$ws.Range("C7").Value = "Ackr4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.1251886666666667
$ws.Range("H7").Value = 0.375566
$ws.Range("I7").Value = 0.03027510681284623
$ws.Range("J7").Value = 0.03027510681284622
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.05974466666666667
$ws.Range("N7").Value = 0.179234
$ws.Range("O7").Value = 0.01783365214557738
$ws.Range("P7").Value = 0.01783365214557738
$ws.Range("Q7").Value = 0.007479355160444445
$ws.Range("R7").Value = 0.06731419644400001
$ws.Range("S7").Value = 0.0005399157235704995
$ws.Range("T7").Value = 0.0005399157235704995
